$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.102.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5252"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2607"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06350"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.645.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5472"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8206"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.122.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.584"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.016"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1239"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.250"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.428"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05895"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.513"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9507"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.411"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5686"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.813"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8475"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.030.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.798.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.890"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.474"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05156"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09686"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
